# Update availability figures on the "Availability" sheet.
# These cells track busy-level per person/hour; change them to 3 (busy)
# to reflect the newly noted meetings / updated availability.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Jan (row 4)
$ws.Range("C4").Value = 3
$ws.Range("P4").Value = 3

# Olivier (row 5)
$ws.Range("C5").Value = 3
$ws.Range("P5").Value = 3

# Pascal (row 6)
$ws.Range("C6").Value = 3
$ws.Range("P6").Value = 3

# Petra (row 7)
$ws.Range("C7").Value = 3
